# "made FP template margins a bit more narrow"
#
# Functional changes applied:
#   1. Selected range moves from T4 to the header row A5:T5
#      (active cell A5, selection A5:T5).
#   2. Left page margin narrowed from 2 cm (0.7874 in) to 1 cm (0.3937 in).
#   3. Print scale bumped from 68% to 72% to compensate for the narrower
#      margin while keeping the sheet fitting the page.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Update the saved selection/active cell for the sheet.
$ws.Range("A5:T5").Select()

# 2. Narrow the left margin to 1 cm (in points: 1/2.54 * 72).
$ws.PageSetup.LeftMargin = 1 / 2.54 * 72

# 3. Bump the print scale up to 72%.
$ws.PageSetup.Zoom = 72
